$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "MODEL_CONDITION" header text to "MODELCONDITION" everywhere it appears.
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")

# The old column A (a running index: 1, 15) is no longer needed; delete it so every
# other column shifts one place to the left (B->A, C->B, D->C, E->D, F->E).
$ws.Range("A1:A3").EntireColumn.Delete()
